$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (last two records removed entirely)
$ws.Rows("4:5").Delete()

# Row 2 becomes the "Instalación" / Chillan / Diwatts record
$ws.Range("B2").Value = 424267
$ws.Range("E2").Value = "GPS, Botón Alámbrico Tablero, Corta Corriente, Sensor Pta, Sensor Pta Adicional, Sensor Temperatura, Sensor Temperatura Adicional"
$ws.Range("F2").Value = "RUTA 5 SUR KM 8 SN"
$ws.Range("G2").Value = "CHILLAN"
$ws.Range("H2").Value = "Región del Ñuble."
$ws.Range("I2").Value = "Pedro Pascal"
$ws.Range("J2").Value = "VE839-POR CONFIRMAR"
$ws.Range("K2").Value = "CIAL_ALIMENTOS"
$ws.Range("L2").Value = "GPS"
$ws.Range("M2").Value = " Botón Alámbrico Tablero"
$ws.Range("N2").Value = " Corta Corriente"
$ws.Range("O2").Value = " Sensor Pta"
$ws.Range("P2").Value = " Sensor Pta Adicional"
$ws.Range("Q2").Value = " Sensor Temperatura"
$ws.Range("R2").Value = " Sensor Temperatura Adicional"

# Row 3 becomes the "Soporte" / Temuco / Juan Perez record
$ws.Range("D3").Value = "Soporte"
$ws.Range("B3").Value = 430947
$ws.Range("E3").Value = "GPS"
$ws.Range("F3").Value = "GUIDO BECK DE RAMBERGA 1884, PADRE DE LAS CASAS "
$ws.Range("G3").Value = "TEMUCO"
$ws.Range("H3").Value = "Región de La Araucanía."
$ws.Range("I3").Value = "Juan Perez"
$ws.Range("J3").Value = "BWYY79"
$ws.Range("K3").Value = "CUENTA_CAROZZIDISTRIBUCION"
$ws.Range("L3").Value = "GPS"
$ws.Range("M3:R3").Clear()

# Restore the selection Excel leaves behind after this edit session
$ws.Range("D5").Select() | Out-Null
